# Updated cryptos list (Price / Volume(1h) columns) to match the latest
# GitHub Actions scrape. Price cells (column D) are forced to Text so
# values such as "4.40" / "0.998" survive verbatim instead of being
# auto-coerced to numbers by the COM Value setter; the style is reset to
# "Normal" afterwards so no stray number-format style id is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.289.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.500.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("E10").Value = "  +2.18%  "
$ws.Range("E11").Value = "  +2.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.092.91"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("E13").Value = "  +1.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000181"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.497.27"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.368.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.65"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.49%  "
$ws.Range("E20").Value = "  -2.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "394.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.83%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.638.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("E26").Value = "  +1.55%  "
$ws.Range("E27").Value = "  +1.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.64%  "
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.80%  "
$ws.Range("E32").Value = "  -6.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.521.26"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("E34").Value = "  +5.15%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.48"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.29%  "
$ws.Range("E37").Value = "  -4.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.89"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.73%  "
$ws.Range("E39").Value = "  -1.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "167.17"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0781"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.809"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.40%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.37"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.66"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.467.73"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.41%  "
$ws.Range("E49").Value = "  -0.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.892"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0258"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.71%  "
